$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value2 = 44568
$ws.Range("J2").Value2 = 500
$ws.Range("K2").Value2 = 15000
$ws.Range("L2").Value2 = 16000
$ws.Range("M2").Value2 = 15500
$ws.Range("P2").Value2 = 861

# Row 4 updates
$ws.Range("D4").Value2 = 44547
$ws.Range("J4").Value2 = 200
$ws.Range("K4").Value2 = 13000
$ws.Range("L4").Value2 = 14000
$ws.Range("M4").Value2 = 13500
$ws.Range("P4").Value2 = 750
